$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '36.690.94'
$ws.Range('D2').ClearFormats()

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.959.94'
$ws.Range('D3').ClearFormats()

$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('E3').ClearFormats()

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('E4').ClearFormats()

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '244.73'
$ws.Range('D5').ClearFormats()

$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.91%  '
$ws.Range('E5').ClearFormats()

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.617'
$ws.Range('D6').ClearFormats()

$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +1.27%  '
$ws.Range('E6').ClearFormats()

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '61.57'
$ws.Range('D7').ClearFormats()

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +8.18%  '
$ws.Range('E7').ClearFormats()

$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E8').ClearFormats()

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.377'
$ws.Range('D9').ClearFormats()

$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +4.91%  '
$ws.Range('E9').ClearFormats()

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0797'
$ws.Range('D10').ClearFormats()

$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -6.36%  '
$ws.Range('E10').ClearFormats()

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.25'
$ws.Range('D12').ClearFormats()

$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +6.22%  '
$ws.Range('E12').ClearFormats()

$ws.Range('B13').NumberFormat = "@"
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('B13').ClearFormats()

$ws.Range('C13').NumberFormat = "@"
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('C13').ClearFormats()

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '22.02'
$ws.Range('D13').ClearFormats()

$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +3.35%  '
$ws.Range('E13').ClearFormats()

$ws.Range('B14').NumberFormat = "@"
$ws.Range('B14').Value = 'Polygon'
$ws.Range('B14').ClearFormats()

$ws.Range('C14').NumberFormat = "@"
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('C14').ClearFormats()

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.835'
$ws.Range('D14').ClearFormats()

$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +3.65%  '
$ws.Range('E14').ClearFormats()

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.236.72'
$ws.Range('D15').ClearFormats()

$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('E15').ClearFormats()

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.29'
$ws.Range('D16').ClearFormats()

$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +3.27%  '
$ws.Range('E16').ClearFormats()

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.979.93'
$ws.Range('D17').ClearFormats()

$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +2.07%  '
$ws.Range('E17').ClearFormats()

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '36.591.03'
$ws.Range('D18').ClearFormats()

$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('E18').ClearFormats()

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '69.89'
$ws.Range('D19').ClearFormats()

$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +1.22%  '
$ws.Range('E19').ClearFormats()

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0854'
$ws.Range('D20').ClearFormats()

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '230.25'
$ws.Range('D21').ClearFormats()

$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +1.42%  '
$ws.Range('E21').ClearFormats()

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.08'
$ws.Range('D22').ClearFormats()

$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +2.23%  '
$ws.Range('E22').ClearFormats()

$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('E23').ClearFormats()

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.46'
$ws.Range('D24').ClearFormats()

$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +5.80%  '
$ws.Range('E24').ClearFormats()

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.34'
$ws.Range('D25').ClearFormats()

$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +3.28%  '
$ws.Range('E25').ClearFormats()

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.141'
$ws.Range('D26').ClearFormats()

$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +5.79%  '
$ws.Range('E26').ClearFormats()

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.20'
$ws.Range('D27').ClearFormats()

$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +0.83%  '
$ws.Range('E27').ClearFormats()

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '160.55'
$ws.Range('D28').ClearFormats()

$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('E28').ClearFormats()

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.42'
$ws.Range('D29').ClearFormats()

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +1.07%  '
$ws.Range('E29').ClearFormats()

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.29'
$ws.Range('D30').ClearFormats()

$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +17.09%  '
$ws.Range('E30').ClearFormats()

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.119'
$ws.Range('D31').ClearFormats()

$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +1.49%  '
$ws.Range('E31').ClearFormats()

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.78'
$ws.Range('D32').ClearFormats()

$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +5.15%  '
$ws.Range('E32').ClearFormats()

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0618'
$ws.Range('D33').ClearFormats()

$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('E33').ClearFormats()

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.47'
$ws.Range('D34').ClearFormats()

$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +8.01%  '
$ws.Range('E34').ClearFormats()

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.53'
$ws.Range('D35').ClearFormats()

$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +15.87%  '
$ws.Range('E35').ClearFormats()

$ws.Range('B36').NumberFormat = "@"
$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('B36').ClearFormats()

$ws.Range('C36').NumberFormat = "@"
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('C36').ClearFormats()

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').ClearFormats()

$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E36').ClearFormats()

$ws.Range('B37').NumberFormat = "@"
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('B37').ClearFormats()

$ws.Range('C37').NumberFormat = "@"
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('C37').ClearFormats()

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.27'
$ws.Range('D37').ClearFormats()

$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +4.78%  '
$ws.Range('E37').ClearFormats()

$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -1.05%  '
$ws.Range('E38').ClearFormats()

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.55'
$ws.Range('D39').ClearFormats()

$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -8.07%  '
$ws.Range('E39').ClearFormats()

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0983'
$ws.Range('D40').ClearFormats()

$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -0.91%  '
$ws.Range('E40').ClearFormats()

$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('E41').ClearFormats()

$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +2.50%  '
$ws.Range('E42').ClearFormats()

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0211'
$ws.Range('D43').ClearFormats()

$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +1.09%  '
$ws.Range('E43').ClearFormats()

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.10'
$ws.Range('D44').ClearFormats()

$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +3.58%  '
$ws.Range('E44').ClearFormats()

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.369.06'
$ws.Range('D45').ClearFormats()

$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +2.34%  '
$ws.Range('E45').ClearFormats()

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '88.73'
$ws.Range('D46').ClearFormats()

$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +3.26%  '
$ws.Range('E46').ClearFormats()

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.03'
$ws.Range('D47').ClearFormats()

$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +1.49%  '
$ws.Range('E47').ClearFormats()

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.15'
$ws.Range('D48').ClearFormats()

$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.73%  '
$ws.Range('E48').ClearFormats()

$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +0.63%  '
$ws.Range('E49').ClearFormats()

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '45.88'
$ws.Range('D50').ClearFormats()

$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +6.10%  '
$ws.Range('E50').ClearFormats()

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.129.50'
$ws.Range('D51').ClearFormats()

$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +0.64%  '
$ws.Range('E51').ClearFormats()
